$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values (price + 1h volume change) scraped on
# Tue Nov 12 07:27:53 UTC 2024. Rows 21/22 (BitcoinCash <-> Uniswap)
# also swapped position in the source ranking.

$ws.Range("D2").Value = '89.672.05'
$ws.Range("E2").Value = '  +10.84%  '
$ws.Range("D3").Value = '3.374.59'
$ws.Range("E3").Value = '  +7.65%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.19'
$ws.Range("E5").Value = '  +6.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '651.99'
$ws.Range("E6").Value = '  +5.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.411'
$ws.Range("E7").Value = '  +47.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").Value = '  +6.89%  '
$ws.Range("D10").Value = '3.370.33'
$ws.Range("E10").Value = '  +7.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.606'
$ws.Range("E11").Value = '  +6.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000290'
$ws.Range("E12").Value = '  +16.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.01'
$ws.Range("E13").Value = '  +18.65%  '
$ws.Range("D15").Value = '4.001.25'
$ws.Range("E15").Value = '  +7.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.60'
$ws.Range("E16").Value = '  +6.92%  '
$ws.Range("D17").Value = '89.554.71'
$ws.Range("E17").Value = '  +10.78%  '
$ws.Range("D18").Value = '3.374.04'
$ws.Range("E18").Value = '  +7.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.85'
$ws.Range("E19").Value = '  +7.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.18'
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.67'
$ws.Range("E21").Value = '  +8.85%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.50'
$ws.Range("E22").Value = '  +7.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.54'
$ws.Range("E23").Value = '  +10.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.54'
$ws.Range("E24").Value = '  +5.78%  '
$ws.Range("E25").Value = '  +10.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.91'
$ws.Range("E26").Value = '  +20.41%  '
$ws.Range("E27").Value = '  +6.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000144'
$ws.Range("E28").Value = '  +20.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '78.90'
$ws.Range("E29").Value = '  +4.60%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.42'
$ws.Range("E32").Value = '  +6.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '598.95'
$ws.Range("E33").Value = '  +9.01%  '
$ws.Range("E34").Value = '  +9.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.11'
$ws.Range("E36").Value = '  +7.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.45'
$ws.Range("E37").Value = '  +25.99%  '
$ws.Range("E38").Value = '  -3.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.50'
$ws.Range("E39").Value = '  +4.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.429'
$ws.Range("E40").Value = '  +6.59%  '
$ws.Range("E41").Value = '  +8.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.23'
$ws.Range("E42").Value = '  +8.38%  '
$ws.Range("E43").Value = '  +5.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.49'
$ws.Range("E45").Value = '  +14.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '158.11'
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '190.38'
$ws.Range("E48").Value = '  +2.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.37'
$ws.Range("E49").Value = '  +3.73%  '
$ws.Range("E50").Value = '  +9.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.674'
$ws.Range("E51").Value = '  +8.88%  '
